# Mark the use cases that are finished as "done" in a new column G on the
# "UC" sheet (this introduces the 52nd shared string "done" and extends the
# sheet's used range from A1:F40 to A1:G40).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UC")

$doneRows = @(2, 3, 10, 12, 13, 14, 23, 29)
foreach ($r in $doneRows) {
    $ws.Cells.Item($r, 7).Value = "done"
}

# Reproduce the author's final on-screen state: scrolled down a bit, with
# G2 (the first newly-entered cell) selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 22 | Out-Null
$excel.ActiveWindow.ScrollColumn = 1 | Out-Null
$ws.Range("G2").Select() | Out-Null
